$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.714.81'
$ws.Range("E2").Value = '  +4.00%  '
$ws.Range("D3").Value = '2.756.68'
$ws.Range("E3").Value = '  +3.95%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '117.96'
$ws.Range("E5").Value = '  +4.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '333.53'
$ws.Range("E6").Value = '  +2.86%  '
$ws.Range("E7").Value = '  +2.53%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +5.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.57'
$ws.Range("E10").Value = '  +3.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.26'
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0832'
$ws.Range("E12").Value = '  +2.04%  '
$ws.Range("E13").Value = '  +3.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.64'
$ws.Range("E14").Value = '  +4.53%  '
$ws.Range("D15").Value = '3.185.98'
$ws.Range("E15").Value = '  +3.86%  '
$ws.Range("D16").Value = '2.753.68'
$ws.Range("E16").Value = '  +4.20%  '
$ws.Range("E17").Value = '  +2.65%  '
$ws.Range("D18").Value = '51.666.60'
$ws.Range("E18").Value = '  +4.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.79'
$ws.Range("E19").Value = '  +6.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.00'
$ws.Range("E20").Value = '  +3.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.88'
$ws.Range("E21").Value = '  +2.56%  '
$ws.Range("E22").Value = '  +1.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '278.34'
$ws.Range("E23").Value = '  +1.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.26'
$ws.Range("E24").Value = '  -0.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.69'
$ws.Range("E25").Value = '  +5.36%  '
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.14'
$ws.Range("E27").Value = '  +0.54%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("E30").Value = '  -0.56%  '
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.141'
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '35.70'
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.54'
$ws.Range("E33").Value = '  +1.59%  '
$ws.Range("E34").Value = '  +3.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0827'
$ws.Range("E35").Value = '  +3.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.49'
$ws.Range("E36").Value = '  +0.54%  '
$ws.Range("E37").Value = '  +3.47%  '
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("E40").Value = '  +5.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '130.10'
$ws.Range("E41").Value = '  +4.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.64'
$ws.Range("E42").Value = '  +6.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0346'
$ws.Range("E43").Value = '  +9.68%  '
$ws.Range("E44").Value = '  +2.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.29'
$ws.Range("E45").Value = '  +3.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.37'
$ws.Range("E46").Value = '  +11.22%  '
$ws.Range("D47").Value = '2.116.46'
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.37'
$ws.Range("E48").Value = '  +3.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.27'
$ws.Range("E49").Value = '  +2.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.66'
$ws.Range("E50").Value = '  +7.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.02'
$ws.Range("E51").Value = '  +0.58%  '
